# Edit: remove the "Fifi" entry row and rename the picture filenames
# (insert an underscore before the dog-name portion of each "path to pic"
# value) as described in the commit "renaming picture file names".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Delete the entire row for "Fifi" (original row 5); rows below shift up.
$ws.Rows(5).Delete()

# After the deletion the data rows (originally 2,3,4,6,7,8) now occupy
# rows 2-7. Update the "path to pic" column (H) with the renamed files.
$ws.Range("H2").Value = "PicturesOrig\FolkeNoertemann_Arlo.jpg"
$ws.Range("H3").Value = "PicturesOrig\AnnieVanderlinck_Moss.jpg"
$ws.Range("H4").Value = "PicturesOrig\AnnieVanderlinck_Tweed.jpg"
$ws.Range("H5").Value = "PicturesOrig\FolkeNoertemannKinloch_Luke.jpg"
$ws.Range("H6").Value = "PicturesOrig\FolkeNoertemannKinloch_Heath.jpg"
$ws.Range("H7").Value = "PicturesOrig\FolkeNoertemann_Joe.jpg"

$ws.Range("H13").Select()

$wb.Save()
Write-Host "Row removed and picture paths renamed."
